$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they keep their
# original text representation (e.g. "597.16") instead of being
# auto-converted to a numeric value by Excel.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated Price / Volume(1h) values from the latest crypto data pull.
$ws.Range("D2").Value = "68.225.36"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "3.698.69"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "597.16"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "167.26"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("D7").Value = "3.697.30"
$ws.Range("E7").Value = "  -3.04%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("D11").Value = "6.27"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("D13").Value = "38.12"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "4.317.13"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").Value = "3.703.67"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("D17").Value = "68.138.62"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "7.26"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "17.24"
$ws.Range("E20").Value = "  +7.28%  "
$ws.Range("D21").Value = "489.40"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Value = "84.57"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("E25").Value = "  +4.20%  "
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("D27").Value = "12.23"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "10.06"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  +3.02%  "
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").Value = "31.41"
$ws.Range("E33").Value = "  -4.65%  "
$ws.Range("D34").Value = "3.841.31"
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "3.644.32"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").Value = "434.51"
$ws.Range("E42").Value = "  -3.60%  "
$ws.Range("D43").Value = "48.60"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").Value = "2.83"
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D48").Value = "40.43"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").Value = "141.41"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").Value = "2.754.80"
$ws.Range("E50").Value = "  -3.14%  "
$ws.Range("E51").Value = "  -0.50%  "
